$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(67, 1).Value = "4wg5vrwdw"
$ws.Cells.Item(67, 2).Value = "wrb1ug0m4"
